$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value2 = 8.296760368086368
$ws.Cells.Item(2, 4).Value2 = 6.782595929419921
$ws.Cells.Item(2, 5).Value2 = 12.55048805425085
$ws.Cells.Item(2, 6).Value2 = 43.75957632028918
$ws.Cells.Item(2, 7).Value2 = 53.0339579370205
$ws.Cells.Item(2, 8).Value2 = 20.45620539983744
$ws.Cells.Item(2, 9).Value2 = 33.79970058325425
$ws.Cells.Item(2, 10).Value2 = 10.46894284426493
$ws.Cells.Item(2, 11).Value2 = 18.76607139671875
$ws.Cells.Item(2, 12).Value2 = 10.21345856491948
$ws.Cells.Item(2, 13).Value2 = 19.97887698724005
$ws.Cells.Item(3, 3).Value2 = 8.269980275549122
$ws.Cells.Item(3, 4).Value2 = 6.76033334489534
$ws.Cells.Item(3, 5).Value2 = 12.56913595964794
$ws.Cells.Item(3, 6).Value2 = 43.91806548420374
$ws.Cells.Item(3, 7).Value2 = 53.25536287100959
$ws.Cells.Item(3, 8).Value2 = 20.54094568574256
$ws.Cells.Item(3, 9).Value2 = 33.92898532588785
$ws.Cells.Item(3, 10).Value2 = 10.4943432626684
$ws.Cells.Item(3, 11).Value2 = 18.39768514395682
$ws.Cells.Item(3, 12).Value2 = 10.23423291124247
$ws.Cells.Item(3, 13).Value2 = 19.83964961765781
$ws.Cells.Item(4, 3).Value2 = 8.25442024269609
$ws.Cells.Item(4, 4).Value2 = 6.746998779562229
$ws.Cells.Item(4, 5).Value2 = 12.58187083143062
$ws.Cells.Item(4, 6).Value2 = 44.02634018647267
$ws.Cells.Item(4, 7).Value2 = 53.40768847762186
$ws.Cells.Item(4, 8).Value2 = 20.59698601571004
$ws.Cells.Item(4, 9).Value2 = 34.01623403274895
$ws.Cells.Item(4, 10).Value2 = 10.5108248446536
$ws.Cells.Item(4, 11).Value2 = 18.17011309525248
$ws.Cells.Item(4, 12).Value2 = 10.24775194412058
$ws.Cells.Item(4, 13).Value2 = 19.75620224175286
$ws.Cells.Item(5, 3).Value2 = 8.248305204450354
$ws.Cells.Item(5, 4).Value2 = 6.741651479197584
$ws.Cells.Item(5, 5).Value2 = 12.58738398285493
$ws.Cells.Item(5, 6).Value2 = 44.07320997382372
$ws.Cells.Item(5, 7).Value2 = 53.47386028299637
$ws.Cells.Item(5, 8).Value2 = 20.62082989921933
$ws.Cells.Item(5, 9).Value2 = 34.05376134042262
$ws.Cells.Item(5, 10).Value2 = 10.5177645036834
$ws.Cells.Item(5, 11).Value2 = 18.07714636376119
$ws.Cells.Item(5, 12).Value2 = 10.2534535524135
$ws.Cells.Item(5, 13).Value2 = 19.72273524048201
$ws.Cells.Item(6, 3).Value2 = 8.247303550830008
$ws.Cells.Item(6, 4).Value2 = 6.740768867910211
$ws.Cells.Item(6, 5).Value2 = 12.58831899485862
$ws.Cells.Item(6, 6).Value2 = 44.08115829511635
$ws.Cells.Item(6, 7).Value2 = 53.48509492310401
$ws.Cells.Item(6, 8).Value2 = 20.62484994274273
$ws.Cells.Item(6, 9).Value2 = 34.0601117062758
$ws.Cells.Item(6, 10).Value2 = 10.51893033156152
$ws.Cells.Item(6, 11).Value2 = 18.06169897151653
$ws.Cells.Item(6, 12).Value2 = 10.25441194105718
$ws.Cells.Item(6, 13).Value2 = 19.71721133141281
$ws.Cells.Item(7, 3).Value2 = 8.254336854059835
$ws.Cells.Item(7, 4).Value2 = 6.746926309960265
$ws.Cells.Item(7, 5).Value2 = 12.58194387290867
$ws.Cells.Item(7, 6).Value2 = 44.02696117921822
$ws.Cells.Item(7, 7).Value2 = 53.40856433053771
$ws.Cells.Item(7, 8).Value2 = 20.59730350701605
$ws.Cells.Item(7, 9).Value2 = 34.01673215903116
$ws.Cells.Item(7, 10).Value2 = 10.51091753049581
$ws.Cells.Item(7, 11).Value2 = 18.16886008347469
$ws.Cells.Item(7, 12).Value2 = 10.2478280579401
$ws.Cells.Item(7, 13).Value2 = 19.75574868031143
$ws.Cells.Item(8, 3).Value2 = 8.287345616998595
$ws.Cells.Item(8, 4).Value2 = 6.774851808928221
$ws.Cells.Item(8, 5).Value2 = 12.55665143802357
$ws.Cells.Item(8, 6).Value2 = 43.81194267999603
$ws.Cells.Item(8, 7).Value2 = 53.10688548157946
$ws.Cells.Item(8, 8).Value2 = 20.48459095611764
$ws.Cells.Item(8, 9).Value2 = 33.84264178143119
$ws.Cells.Item(8, 10).Value2 = 10.47751744928301
$ws.Cells.Item(8, 11).Value2 = 18.63940837106576
$ws.Cells.Item(8, 12).Value2 = 10.22046339868336
$ws.Cells.Item(8, 13).Value2 = 19.9304639607819
$ws.Cells.Item(9, 3).Value2 = 8.358910816056033
$ws.Cells.Item(9, 4).Value2 = 6.832159958100127
$ws.Cells.Item(9, 5).Value2 = 12.51722817977225
$ws.Cells.Item(9, 6).Value2 = 43.47768021475938
$ws.Cells.Item(9, 7).Value2 = 52.64621612165902
$ws.Cells.Item(9, 8).Value2 = 20.29542975427037
$ws.Cells.Item(9, 9).Value2 = 33.56391994483657
$ws.Cells.Item(9, 10).Value2 = 10.41902043085255
$ws.Cells.Item(9, 11).Value2 = 19.54618217353678
$ws.Cells.Item(9, 12).Value2 = 10.17283657850683
$ws.Cells.Item(9, 13).Value2 = 20.28809001822064
$ws.Cells.Item(10, 3).Value2 = 8.4154174721325
$ws.Cells.Item(10, 4).Value2 = 6.875667591541499
$ws.Cells.Item(10, 5).Value2 = 12.49443911907041
$ws.Cells.Item(10, 6).Value2 = 43.28591905389494
$ws.Cells.Item(10, 7).Value2 = 52.38881911095737
$ws.Cells.Item(10, 8).Value2 = 20.17595558395969
$ws.Cells.Item(10, 9).Value2 = 33.39767728172136
$ws.Cells.Item(10, 10).Value2 = 10.38027320648134
$ws.Cells.Item(10, 11).Value2 = 20.19624846464786
$ws.Cells.Item(10, 12).Value2 = 10.14149264989818
$ws.Cells.Item(10, 13).Value2 = 20.55838528157347
$ws.Cells.Item(11, 3).Value2 = 8.441922424693152
$ws.Cells.Item(11, 4).Value2 = 6.895732279460462
$ws.Cells.Item(11, 5).Value2 = 12.48540654301863
$ws.Cells.Item(11, 6).Value2 = 43.21047962703716
$ws.Cells.Item(11, 7).Value2 = 52.28957508316002
$ws.Cells.Item(11, 8).Value2 = 20.12585516232956
$ws.Cells.Item(11, 9).Value2 = 33.33048500361585
$ws.Cells.Item(11, 10).Value2 = 10.36355684615439
$ws.Cells.Item(11, 11).Value2 = 20.48721995537754
$ws.Cells.Item(11, 12).Value2 = 10.12801878889649
$ws.Cells.Item(11, 13).Value2 = 20.68263750976348
$ws.Cells.Item(12, 3).Value2 = 8.452069314304664
$ws.Cells.Item(12, 4).Value2 = 6.903366599531598
$ws.Cells.Item(12, 5).Value2 = 12.48217747707006
$ws.Cells.Item(12, 6).Value2 = 43.18361716979923
$ws.Cells.Item(12, 7).Value2 = 52.25457982620726
$ws.Cells.Item(12, 8).Value2 = 20.10749593640582
$ws.Cells.Item(12, 9).Value2 = 33.30625907865601
$ws.Cells.Item(12, 10).Value2 = 10.35735705325512
$ws.Cells.Item(12, 11).Value2 = 20.59662353019956
$ws.Cells.Item(12, 12).Value2 = 10.12302891012766
$ws.Cells.Item(12, 13).Value2 = 20.72984455910571
$ws.Cells.Item(13, 3).Value2 = 8.449879186600644
$ws.Cells.Item(13, 4).Value2 = 6.901720847368265
$ws.Cells.Item(13, 5).Value2 = 12.48286441078975
$ws.Cells.Item(13, 6).Value2 = 43.1893265210166
$ws.Cells.Item(13, 7).Value2 = 52.26200134720166
$ws.Cells.Item(13, 8).Value2 = 20.11142264585393
$ws.Cells.Item(13, 9).Value2 = 33.31142228910267
$ws.Cells.Item(13, 10).Value2 = 10.35868650133289
$ws.Cells.Item(13, 11).Value2 = 20.57309775142364
$ws.Cells.Item(13, 12).Value2 = 10.12409857881717
$ws.Cells.Item(13, 13).Value2 = 20.71967125186666
$ws.Cells.Item(14, 3).Value2 = 8.442755038996879
$ws.Cells.Item(14, 4).Value2 = 6.896359650099545
$ws.Cells.Item(14, 5).Value2 = 12.48513705327086
$ws.Cells.Item(14, 6).Value2 = 43.20823541992158
$ws.Cells.Item(14, 7).Value2 = 52.28664405819503
$ws.Cells.Item(14, 8).Value2 = 20.12433244809233
$ws.Cells.Item(14, 9).Value2 = 33.32846747247963
$ws.Cells.Item(14, 10).Value2 = 10.36304417640981
$ws.Cells.Item(14, 11).Value2 = 20.49623684824575
$ws.Cells.Item(14, 12).Value2 = 10.12760601887693
$ws.Cells.Item(14, 13).Value2 = 20.68651832272116
$ws.Cells.Item(15, 3).Value2 = 8.4384054734663
$ws.Cells.Item(15, 4).Value2 = 6.893080394785557
$ws.Cells.Item(15, 5).Value2 = 12.48655402011524
$ws.Cells.Item(15, 6).Value2 = 43.22003994112888
$ws.Cells.Item(15, 7).Value2 = 52.30207581834203
$ws.Cells.Item(15, 8).Value2 = 20.13231991987587
$ws.Cells.Item(15, 9).Value2 = 33.33906697100705
$ws.Cells.Item(15, 10).Value2 = 10.36573033571112
$ws.Cells.Item(15, 11).Value2 = 20.44905279683704
$ws.Cells.Item(15, 12).Value2 = 10.12976904904753
$ws.Cells.Item(15, 13).Value2 = 20.66623052781016
$ws.Cells.Item(16, 3).Value2 = 8.413700983284947
$ws.Cells.Item(16, 4).Value2 = 6.874361553906446
$ws.Cells.Item(16, 5).Value2 = 12.49505622648869
$ws.Cells.Item(16, 6).Value2 = 43.29108728677842
$ws.Cells.Item(16, 7).Value2 = 52.39566579388947
$ws.Cells.Item(16, 8).Value2 = 20.17931541102322
$ws.Cells.Item(16, 9).Value2 = 33.40223863225375
$ws.Cells.Item(16, 10).Value2 = 10.38138392587952
$ws.Cells.Item(16, 11).Value2 = 20.17712872356746
$ws.Cells.Item(16, 12).Value2 = 10.14238894884574
$ws.Cells.Item(16, 13).Value2 = 20.55028852458104
$ws.Cells.Item(17, 3).Value2 = 8.39874693230607
$ws.Cells.Item(17, 4).Value2 = 6.862945930041042
$ws.Cells.Item(17, 5).Value2 = 12.50061345009838
$ws.Cells.Item(17, 6).Value2 = 43.33769963565626
$ws.Cells.Item(17, 7).Value2 = 52.45766574224334
$ws.Cells.Item(17, 8).Value2 = 20.20923510342599
$ws.Cells.Item(17, 9).Value2 = 33.44315627996458
$ws.Cells.Item(17, 10).Value2 = 10.39121958436029
$ws.Cells.Item(17, 11).Value2 = 20.00902439153694
$ws.Cells.Item(17, 12).Value2 = 10.15033149340811
$ws.Cells.Item(17, 13).Value2 = 20.47947145921181
$ws.Cells.Item(18, 3).Value2 = 8.390221344313105
$ws.Cells.Item(18, 4).Value2 = 6.856405827880146
$ws.Cells.Item(18, 5).Value2 = 12.50393543756892
$ws.Cells.Item(18, 6).Value2 = 43.36561939556209
$ws.Cells.Item(18, 7).Value2 = 52.49500519622019
$ws.Cells.Item(18, 8).Value2 = 20.22684394865794
$ws.Cells.Item(18, 9).Value2 = 33.4674843845488
$ws.Cells.Item(18, 10).Value2 = 10.39696247411041
$ws.Cells.Item(18, 11).Value2 = 19.91189393342217
$ws.Cells.Item(18, 12).Value2 = 10.15497371659467
$ws.Cells.Item(18, 13).Value2 = 20.4388629548421
$ws.Cells.Item(19, 3).Value2 = 8.387347857944901
$ws.Cells.Item(19, 4).Value2 = 6.854195993585738
$ws.Cells.Item(19, 5).Value2 = 12.50508179548757
$ws.Cells.Item(19, 6).Value2 = 43.37526283866339
$ws.Cells.Item(19, 7).Value2 = 52.50793539526762
$ws.Cells.Item(19, 8).Value2 = 20.23287462025684
$ws.Cells.Item(19, 9).Value2 = 33.47585754413651
$ws.Cells.Item(19, 10).Value2 = 10.39892164929532
$ws.Cells.Item(19, 11).Value2 = 19.8789345336287
$ws.Cells.Item(19, 12).Value2 = 10.15655819667859
$ws.Cells.Item(19, 13).Value2 = 20.42513577073722
$ws.Cells.Item(20, 3).Value2 = 8.400331029212476
$ws.Cells.Item(20, 4).Value2 = 6.864158487430339
$ws.Cells.Item(20, 5).Value2 = 12.50000887654664
$ws.Cells.Item(20, 6).Value2 = 43.33262276943049
$ws.Cells.Item(20, 7).Value2 = 52.45089185317041
$ws.Cells.Item(20, 8).Value2 = 20.20600870781673
$ws.Cells.Item(20, 9).Value2 = 33.43871837480479
$ws.Cells.Item(20, 10).Value2 = 10.39016369802821
$ws.Cells.Item(20, 11).Value2 = 20.02696578870986
$ws.Cells.Item(20, 12).Value2 = 10.14947835279148
$ws.Cells.Item(20, 13).Value2 = 20.48699748583894
$ws.Cells.Item(21, 3).Value2 = 8.444844627674541
$ws.Cells.Item(21, 4).Value2 = 6.897933403592363
$ws.Cells.Item(21, 5).Value2 = 12.48446433324463
$ws.Cells.Item(21, 6).Value2 = 43.20263508080304
$ws.Cells.Item(21, 7).Value2 = 52.27933555100746
$ws.Cells.Item(21, 8).Value2 = 20.12052388347815
$ws.Cells.Item(21, 9).Value2 = 33.32342777850317
$ws.Cells.Item(21, 10).Value2 = 10.36176068860069
$ws.Cells.Item(21, 11).Value2 = 20.51883470950995
$ws.Cells.Item(21, 12).Value2 = 10.12657275232738
$ws.Cells.Item(21, 13).Value2 = 20.69625216770615
$ws.Cells.Item(22, 3).Value2 = 8.47457583397915
$ws.Cells.Item(22, 4).Value2 = 6.920217194747111
$ws.Cells.Item(22, 5).Value2 = 12.47542031291967
$ws.Cells.Item(22, 6).Value2 = 43.12762142561498
$ws.Cells.Item(22, 7).Value2 = 52.18229711903683
$ws.Cells.Item(22, 8).Value2 = 20.06822680250245
$ws.Cells.Item(22, 9).Value2 = 33.25518265065667
$ws.Cells.Item(22, 10).Value2 = 10.34395710574109
$ws.Cells.Item(22, 11).Value2 = 20.83570577057497
$ws.Cells.Item(22, 12).Value2 = 10.1122574482498
$ws.Cells.Item(22, 13).Value2 = 20.83390390642993
$ws.Cells.Item(23, 3).Value2 = 8.458650937115605
$ws.Cells.Item(23, 4).Value2 = 6.908305701457011
$ws.Cells.Item(23, 5).Value2 = 12.48014539960588
$ws.Cells.Item(23, 6).Value2 = 43.16674520561512
$ws.Cells.Item(23, 7).Value2 = 52.23270187524529
$ws.Cells.Item(23, 8).Value2 = 20.09581130937882
$ws.Cells.Item(23, 9).Value2 = 33.29095448192651
$ws.Cells.Item(23, 10).Value2 = 10.35338989560926
$ws.Cells.Item(23, 11).Value2 = 20.66703699436557
$ws.Cells.Item(23, 12).Value2 = 10.11983802657639
$ws.Cells.Item(23, 13).Value2 = 20.76036511913896
$ws.Cells.Item(24, 3).Value2 = 8.399614635328328
$ws.Cells.Item(24, 4).Value2 = 6.863610218728482
$ws.Cells.Item(24, 5).Value2 = 12.50028180841942
$ws.Cells.Item(24, 6).Value2 = 43.33491452716861
$ws.Cells.Item(24, 7).Value2 = 52.45394905037241
$ws.Cells.Item(24, 8).Value2 = 20.20746609162731
$ws.Cells.Item(24, 9).Value2 = 33.44072224763914
$ws.Cells.Item(24, 10).Value2 = 10.39064078933274
$ws.Cells.Item(24, 11).Value2 = 20.01885598627384
$ws.Cells.Item(24, 12).Value2 = 10.14986382105113
$ws.Cells.Item(24, 13).Value2 = 20.48359463925753
$ws.Cells.Item(25, 3).Value2 = 8.338842798077113
$ws.Cells.Item(25, 4).Value2 = 6.81639912585536
$ws.Cells.Item(25, 5).Value2 = 12.52680671834074
$ws.Cells.Item(25, 6).Value2 = 43.55869379624154
$ws.Cells.Item(25, 7).Value2 = 52.75669220672186
$ws.Cells.Item(25, 8).Value2 = 20.34318354462696
$ws.Cells.Item(25, 9).Value2 = 33.63257796306301
$ws.Cells.Item(25, 10).Value2 = 10.4340998362014
$ws.Cells.Item(25, 11).Value2 = 19.30326564752928
$ws.Cells.Item(25, 12).Value2 = 10.18507808106878
$ws.Cells.Item(25, 13).Value2 = 20.18989913149093
